$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 112094769
$ws.Range("B2").Value = 78739
$ws.Range("E2").Value = 6461
$ws.Range("F2").Value = "Norrlandslav"
$ws.Range("G2").Value = "Nephroma arcticum"
$ws.Range("H2").Value = "(L.) Torss."

# Row 3 updates
$ws.Range("A3").Value = 112094770
$ws.Range("B3").Value = 94048
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 2869
$ws.Range("F3").Value = "Bollvitmossa"
$ws.Range("G3").Value = "Sphagnum wulfianum"
$ws.Range("H3").Value = "Girg."
$ws.Range("S3").Value = 50

# Row 4 updates
$ws.Range("A4").Value = 112094771
$ws.Range("B4").Value = 77650
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("S4").Value = 10
